# Fruta / hortaliza, semanal
# Insert a new weekly price-report row at row 163 (pushing the existing
# rows 163-169 down to 164-170), and populate it with the new week's
# data for Pomelo "Start Ruby" at Feria Lagunitas de Puerto Montt.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Shift rows 163..169 down to 164..170, leaving a blank row 163.
$ws.Rows.Item(163).Insert()

# Fill in the new row 163 with this week's values.
$ws.Cells.Item(163, 1).Value = 4
$ws.Cells.Item(163, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(163, 3).Value = "Los Lagos"
$ws.Cells.Item(163, 4).Value = 44509
$ws.Cells.Item(163, 5).Value = 10
$ws.Cells.Item(163, 6).Value = "Fruta"
$ws.Cells.Item(163, 7).Value = 100102
$ws.Cells.Item(163, 8).Value = "Cítricos"
$ws.Cells.Item(163, 9).Value = 100102006
$ws.Cells.Item(163, 10).Value = "Pomelo"
$ws.Cells.Item(163, 11).Value = "Start Ruby"
$ws.Cells.Item(163, 12).Value = "Primera"
$ws.Cells.Item(163, 13).Value = 300
$ws.Cells.Item(163, 14).Value = 11000
$ws.Cells.Item(163, 15).Value = 12000
$ws.Cells.Item(163, 16).Value = 11500
$ws.Cells.Item(163, 17).Value = "$/caja 14 kilos empedrada"
$ws.Cells.Item(163, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(163, 19).Value = 821
$ws.Cells.Item(163, 20).Value = 14
